$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B ("Day 10 Solution") corrections for data rows 2-21 ---
$bValues = @{
    2  = 182
    3  = 182
    4  = 180
    5  = 182
    6  = 185
    7  = 181
    8  = 183
    21 = 10
}
foreach ($row in $bValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $bValues[$row]
}

# --- Column D ("Total") now mirrors column C ("Day 10 Night") for data rows 2-21 ---
$dValues = @{
    2  = 388.29
    3  = 392.0411
    4  = 385.2151
    5  = 384.4421
    6  = 390.0327
    7  = 383.1463
    8  = 388.8807
    9  = 401.2084
    10 = 7.0614
    11 = 7.0614
    12 = 7.0614
    13 = 7.0614
    14 = 7.0614
    15 = 7.0614
    16 = 5.6491
    17 = 5.6491
    18 = 5.6491
    19 = 409.8895
    20 = 5.6491
    21 = 5.6491
}
foreach ($row in $dValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $dValues[$row]
}

# --- Row 22 ("Total Solution" row): B22 becomes the solution-column total, D22 mirrors B22 ---
$ws.Range("A22").Value = "Total Solution"
$ws.Range("B22").Value = 1763
$ws.Range("D22").Value = 1763

# --- New row 23 ("Total Night" row): sums the night column ---
$ws.Range("A23").Value = "Total Night"
$ws.Range("B23").Value = 1763
$ws.Range("C23").Value = 3593.7598000000007
$ws.Range("D23").Value = 3593.7598000000007

# --- Header rename: D1 "Total" -> "Total Night" (same shared text as A23) ---
$ws.Range("D1").Value = "Total Night"
